# Edit script for horarios-141-completo.xlsx
# Commit: "141: 30/12 19:38 LP1912+6203+6173"
# Adds new scraped bus-arrival rows to sheet "LP1912" (rows 386-408)
# and sheet "6203-6173" (rows 54-56); refreshes the "Ultima actualizacion"
# timestamp and "Total filas" counters on the affected sheets.

$wb = $excel.ActiveWorkbook

$sheet1Data = @(
  ,@(386, "16:38:24", "16:42", "16_P MOR-SANTA ANA", 4, "LP1912", "30/12/2025")
  ,@(387, "16:38:24", "16:48", "15_ABASTO", 10, "LP1912", "30/12/2025")
  ,@(388, "16:38:24", "16:56", "17_179 Y 38", 18, "LP1912", "30/12/2025")
  ,@(389, "16:38:24", "16:57", "10_OLMOS", 19, "LP1912", "30/12/2025")
  ,@(390, "16:38:24", "17:00", "16_SANTA ANA", 22, "LP1912", "30/12/2025")
  ,@(391, "16:38:24", "17:04", "11_ETCHEVERRY", 26, "LP1912", "30/12/2025")
  ,@(392, "16:38:24", "17:04", "23_HERNANDEZ", 26, "LP1912", "30/12/2025")
  ,@(393, "16:38:24", "17:10", "10_OLMOS", 32, "LP1912", "30/12/2025")
  ,@(394, "16:38:24", "17:21", "26_HERNANDEZ", 43, "LP1912", "30/12/2025")
  ,@(395, "16:38:24", "17:22", "10_OLMOS", 44, "LP1912", "30/12/2025")
  ,@(396, "16:38:24", "17:24", "84_COLONIA URQUIZA-ESC 49", 46, "LP1912", "30/12/2025")
  ,@(397, "16:38:24", "17:28", "14_ABASTO", 50, "LP1912", "30/12/2025")
  ,@(398, "16:38:24", "17:31", "15_ABASTO", 53, "LP1912", "30/12/2025")
  ,@(399, "16:38:24", "17:37", "27_EL RETIRO", 59, "LP1912", "30/12/2025")
  ,@(400, "16:38:24", "17:38", "17_ROMERO", 60, "LP1912", "30/12/2025")
  ,@(401, "16:38:24", "17:40", "16_SANTA ANA", 62, "LP1912", "30/12/2025")
  ,@(402, "16:38:24", "17:45", "15_ABASTO", 67, "LP1912", "30/12/2025")
  ,@(403, "16:38:24", "17:50", "16_P MOR-167 Y 521", 72, "LP1912", "30/12/2025")
  ,@(404, "16:38:24", "17:52", "81_EL PELIGRO", 74, "LP1912", "30/12/2025")
  ,@(405, "16:38:24", "17:56", "23_HERNANDEZ", 78, "LP1912", "30/12/2025")
  ,@(406, "16:38:24", "18:04", "17_ROMERO", 86, "LP1912", "30/12/2025")
  ,@(407, "16:38:24", "18:06", "23_HERNANDEZ", 88, "LP1912", "30/12/2025")
  ,@(408, "16:38:24", "18:15", "16_SANTA ANA", 97, "LP1912", "30/12/2025")
)

$sheet3Data = @(
  ,@(54, "30/12/2025", "16:38:35", "16:53", "215B_LP-P MOR-40 Y 115", 15, "L6173")
  ,@(55, "30/12/2025", "16:38:35", "17:26", "215A_LA PLATA", 48, "L6173")
  ,@(56, "30/12/2025", "16:38:30", "18:04", "215C_LA PLATA", 86, "L6203")
)


$newTimestamp = "30/12/2025 16:38:35"

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

foreach ($r in $sheet1Data) {
  $rn = $r[0]
  $ws1.Cells.Item($rn, 1).Value = ""
  $ws1.Cells.Item($rn, 2).Value = $r[1]
  $ws1.Cells.Item($rn, 3).Value = $r[2]
  $ws1.Cells.Item($rn, 4).Value = $r[3]
  $ws1.Cells.Item($rn, 5).Value = $r[4]
  $ws1.Cells.Item($rn, 6).Value = $r[5]
  $ws1.Cells.Item($rn, 7).Value = $r[6]
}

$ws1.Cells.Item(2, 1).Value = "Última actualización: " + $newTimestamp
$ws1.Cells.Item(3, 1).Value = "Total filas: 407"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 (only the timestamp banner changes)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: " + $newTimestamp

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

foreach ($r in $sheet3Data) {
  $rn = $r[0]
  $ws3.Cells.Item($rn, 1).Value = ""
  $ws3.Cells.Item($rn, 2).Value = $r[1]
  $ws3.Cells.Item($rn, 3).Value = $r[2]
  $ws3.Cells.Item($rn, 4).Value = $r[3]
  $ws3.Cells.Item($rn, 5).Value = $r[4]
  $ws3.Cells.Item($rn, 6).Value = $r[5]
  $ws3.Cells.Item($rn, 7).Value = $r[6]
}

$ws3.Cells.Item(2, 1).Value = "Última actualización: " + $newTimestamp
$ws3.Cells.Item(3, 1).Value = "Total filas: 55"
